$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly refresh prepends two new price records (new rows 186 & 187)
# ahead of the existing history, pushing everything from the old row 186
# down by two rows (old 186 -> new 188, ..., old 247 -> new 249).
$ws.Rows("186:187").Insert()

# New row 186: Ají Cristal, Primera, $/saco 25 kilos, Región del Maule
$ws.Range("A186").Value = 4
$ws.Range("B186").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C186").Value = "Los Lagos"
$ws.Range("D186").Value = 44663
$ws.Range("E186").Value = 10
$ws.Range("F186").Value = 100112021
$ws.Range("G186").Value = "Ají"
$ws.Range("H186").Value = "Cristal"
$ws.Range("I186").Value = "Primera"
$ws.Range("J186").Value = 60
$ws.Range("K186").Value = 20000
$ws.Range("L186").Value = 20000
$ws.Range("M186").Value = 20000
$ws.Range("N186").Value = "$/saco 25 kilos"
$ws.Range("O186").Value = "Región del Maule"
$ws.Range("P186").Value = 800
$ws.Range("Q186").Value = 25
$ws.Range("R186").Value = "Hortaliza"

# New row 187: Ají Inferno, Primera, $/caja 15 kilos, Provincia de Quillota
$ws.Range("A187").Value = 4
$ws.Range("B187").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C187").Value = "Los Lagos"
$ws.Range("D187").Value = 44663
$ws.Range("E187").Value = 10
$ws.Range("F187").Value = 100112021
$ws.Range("G187").Value = "Ají"
$ws.Range("H187").Value = "Inferno"
$ws.Range("I187").Value = "Primera"
$ws.Range("J187").Value = 120
$ws.Range("K187").Value = 30000
$ws.Range("L187").Value = 30000
$ws.Range("M187").Value = 30000
$ws.Range("N187").Value = "$/caja 15 kilos"
$ws.Range("O187").Value = "Provincia de Quillota"
$ws.Range("P187").Value = 2000
$ws.Range("Q187").Value = 15
$ws.Range("R187").Value = "Hortaliza"

# Apply the date-cell number format (yyyy-mm-dd hh:mm:ss) used by the rest
# of column D to the two freshly inserted date cells, matching D188's style.
$ws.Range("D186:D187").NumberFormat = $ws.Range("D188").NumberFormat
